$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Fix the mislabeled "Signalwandlung Ausgabebereich" headers -
#    they actually describe the INPUT range into the signal
#    conversion (resistance), not the output range.
# -----------------------------------------------------------------
$ws.Range("AA1").Value = "Signalwandlung Eingabebereich von"
$ws.Range("AB1").Value = "Signalwandlung Eingabebereich bis"
$ws.Range("AC1").Value = "Signalwandlung Eingabebereich unit"

# -----------------------------------------------------------------
# 2. Reorder / fix the affected dataset in row 8 (valve whose
#    Signalwandlung is a resistance-based sensor).
#
#    Previously:
#      T8/U8   = raw sensor output range (plain numbers)
#      V8      = unit "A"
#      W8      = "Widerstand"
#      Y8      = "Ohm"
#      AA8/AB8 = formulas T8*X8 / U8*X8 (computed output range)
#      AC8     = unit "V"
#
#    Now (the raw values actually belong to the INPUT of the signal
#    conversion, and the computed, resistance-scaled values are the
#    real sensor OUTPUT range):
#      T8/U8   = formulas AA8*X8 / AB8*X8 (computed output range)
#      V8      = unit "V"
#      W8      = "Widerstand"
#      Y8      = "Ohm"
#      AA8/AB8 = raw input values (plain numbers)
#      AC8     = unit "A"
# -----------------------------------------------------------------

# Remember the original raw values/labels before overwriting anything.
$rawVon = $ws.Range("T8").Value()
$rawBis = $ws.Range("U8").Value()

# Move the raw input values into their corrected home (AA/AB) as plain
# numbers.
$ws.Range("AA8").Value = $rawVon
$ws.Range("AB8").Value = $rawBis
$ws.Range("AC8").Value = "A"

# The computed output range now lives in T/U, driven off the
# (relocated) input values in AA/AB.
$ws.Range("T8").Formula = "=AA8*X8"
$ws.Range("U8").Formula = "=AB8*X8"
$ws.Range("V8").Value = "V"

# Downstream formulas that referenced the old T8/U8 input values must
# now reference AA8/AB8 instead.
$ws.Range("AD8").Formula = "=(Q8-P8)/(AB8-AA8)"
$ws.Range("AE8").Formula = "=Q8-(AD8*AB8)"

# -----------------------------------------------------------------
# 3. Selection / view bookkeeping to match the saved state.
# -----------------------------------------------------------------
$ws.Range("AB3").Select()
